{"js": "// Apply the dated worksheet update: refresh the header date and all\n// three-digit \u00f7 one-digit division problems/answers in the table.\nconst replacements = [\n  [\"2025-11-04 Tuesday\", \"2025-11-05 Wednesday\"],\n  [\"303\u00f73=101, 0\", \"528\u00f75=105, 3\"],\n  [\"946\u00f74=236, 2\", \"421\u00f77=60, 1\"],\n  [\"988\u00f73=329, 1\", \"840\u00f75=168, 0\"],\n  [\"773\u00f74=193, 1\", \"255\u00f73=85, 0\"],\n  [\"520\u00f73=173, 1\", \"541\u00f73=180, 1\"],\n  [\"981\u00f72=490, 1\", \"229\u00f79=25, 4\"],\n  [\"626\u00f75=125, 1\", \"688\u00f74=172, 0\"],\n  [\"800\u00f72=400, 0\", \"242\u00f73=80, 2\"],\n  [\"121\u00f77=17, 2\", \"700\u00f75=140, 0\"],\n  [\"467\u00f72=233, 1\", \"778\u00f78=97, 2\"],\n  [\"290\u00f77=41, 3\", \"702\u00f76=117, 0\"],\n  [\"103\u00f77=14, 5\", \"102\u00f78=12, 6\"],\n  [\"812\u00f73=270, 2\", \"289\u00f72=144, 1\"],\n  [\"812\u00f79=90, 2\", \"404\u00f75=80, 4\"],\n  [\"375\u00f75=75, 0\", \"947\u00f78=118, 3\"],\n  [\"178\u00f74=44, 2\", \"214\u00f78=26, 6\"],\n  [\"808\u00f74=202, 0\", \"491\u00f73=163, 2\"],\n  [\"444\u00f79=49, 3\", \"685\u00f75=137, 0\"],\n  [\"118\u00f75=23, 3\", \"960\u00f76=160, 0\"],\n  [\"575\u00f79=63, 8\", \"715\u00f76=119, 1\"],\n  [\"439\u00f75=87, 4\", \"868\u00f72=434, 0\"],\n  [\"401\u00f78=50, 1\", \"883\u00f74=220, 3\"],\n  [\"900\u00f76=150, 0\", \"438\u00f75=87, 3\"],\n  [\"997\u00f78=124, 5\", \"223\u00f76=37, 1\"],\n  [\"989\u00f75=197, 4\", \"713\u00f74=178, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the dated worksheet update: refresh the header date and all\n# three-digit \u00f7 one-digit division problems/answers in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-11-04 Tuesday\", \"2025-11-05 Wednesday\"),\n    @(\"303\u00f73=101, 0\", \"528\u00f75=105, 3\"),\n    @(\"946\u00f74=236, 2\", \"421\u00f77=60, 1\"),\n    @(\"988\u00f73=329, 1\", \"840\u00f75=168, 0\"),\n    @(\"773\u00f74=193, 1\", \"255\u00f73=85, 0\"),\n    @(\"520\u00f73=173, 1\", \"541\u00f73=180, 1\"),\n    @(\"981\u00f72=490, 1\", \"229\u00f79=25, 4\"),\n    @(\"626\u00f75=125, 1\", \"688\u00f74=172, 0\"),\n    @(\"800\u00f72=400, 0\", \"242\u00f73=80, 2\"),\n    @(\"121\u00f77=17, 2\", \"700\u00f75=140, 0\"),\n    @(\"467\u00f72=233, 1\", \"778\u00f78=97, 2\"),\n    @(\"290\u00f77=41, 3\", \"702\u00f76=117, 0\"),\n    @(\"103\u00f77=14, 5\", \"102\u00f78=12, 6\"),\n    @(\"812\u00f73=270, 2\", \"289\u00f72=144, 1\"),\n    @(\"812\u00f79=90, 2\", \"404\u00f75=80, 4\"),\n    @(\"375\u00f75=75, 0\", \"947\u00f78=118, 3\"),\n    @(\"178\u00f74=44, 2\", \"214\u00f78=26, 6\"),\n    @(\"808\u00f74=202, 0\", \"491\u00f73=163, 2\"),\n    @(\"444\u00f79=49, 3\", \"685\u00f75=137, 0\"),\n    @(\"118\u00f75=23, 3\", \"960\u00f76=160, 0\"),\n    @(\"575\u00f79=63, 8\", \"715\u00f76=119, 1\"),\n    @(\"439\u00f75=87, 4\", \"868\u00f72=434, 0\"),\n    @(\"401\u00f78=50, 1\", \"883\u00f74=220, 3\"),\n    @(\"900\u00f76=150, 0\", \"438\u00f75=87, 3\"),\n    @(\"997\u00f78=124, 5\", \"223\u00f76=37, 1\"),\n    @(\"989\u00f75=197, 4\", \"713\u00f74=178, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
